{"js": "// Replace each two-digit-by-two-digit multiplication problem's text\n// with its updated equivalent, matching the unified diff exactly.\n// Each <w:t> run in the table holds a unique 'A\u00d7B=C' string, so an\n// exact, case-sensitive whole-string search/replace is unambiguous.\nconst replacements = [\n  [\"23\u00d746=1058\", \"96\u00d758=5568\"],\n  [\"61\u00d744=2684\", \"53\u00d779=4187\"],\n  [\"28\u00d731=868\", \"90\u00d780=7200\"],\n  [\"98\u00d728=2744\", \"89\u00d722=1958\"],\n  [\"35\u00d776=2660\", \"22\u00d739=858\"],\n  [\"12\u00d756=672\", \"33\u00d754=1782\"],\n  [\"26\u00d787=2262\", \"83\u00d745=3735\"],\n  [\"68\u00d711=748\", \"85\u00d725=2125\"],\n  [\"76\u00d723=1748\", \"80\u00d761=4880\"],\n  [\"76\u00d746=3496\", \"30\u00d747=1410\"],\n  [\"73\u00d791=6643\", \"24\u00d744=1056\"],\n  [\"36\u00d713=468\", \"45\u00d794=4230\"],\n  [\"38\u00d754=2052\", \"92\u00d795=8740\"],\n  [\"98\u00d714=1372\", \"92\u00d770=6440\"],\n  [\"12\u00d742=504\", \"95\u00d764=6080\"],\n  [\"67\u00d755=3685\", \"92\u00d768=6256\"],\n  [\"26\u00d776=1976\", \"21\u00d718=378\"],\n  [\"55\u00d713=715\", \"26\u00d730=780\"],\n  [\"85\u00d727=2295\", \"73\u00d770=5110\"],\n  [\"76\u00d736=2736\", \"61\u00d719=1159\"],\n  [\"47\u00d754=2538\", \"98\u00d719=1862\"],\n  [\"56\u00d777=4312\", \"36\u00d717=612\"],\n  [\"57\u00d738=2166\", \"70\u00d799=6930\"],\n  [\"82\u00d729=2378\", \"73\u00d719=1387\"],\n  [\"28\u00d776=2128\", \"48\u00d731=1488\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each two-digit-by-two-digit multiplication problem's text\n# to its new value, matching the unified diff exactly. Every run in the\n# table holds a unique 'A*B=C' string, so a literal Find/Replace over the\n# whole document body is unambiguous (one hit per pair) and preserves the\n# existing run formatting (font / size) since only Find.Text changes.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"23\u00d746=1058\", \"96\u00d758=5568\"),\n    @(\"61\u00d744=2684\", \"53\u00d779=4187\"),\n    @(\"28\u00d731=868\", \"90\u00d780=7200\"),\n    @(\"98\u00d728=2744\", \"89\u00d722=1958\"),\n    @(\"35\u00d776=2660\", \"22\u00d739=858\"),\n    @(\"12\u00d756=672\", \"33\u00d754=1782\"),\n    @(\"26\u00d787=2262\", \"83\u00d745=3735\"),\n    @(\"68\u00d711=748\", \"85\u00d725=2125\"),\n    @(\"76\u00d723=1748\", \"80\u00d761=4880\"),\n    @(\"76\u00d746=3496\", \"30\u00d747=1410\"),\n    @(\"73\u00d791=6643\", \"24\u00d744=1056\"),\n    @(\"36\u00d713=468\", \"45\u00d794=4230\"),\n    @(\"38\u00d754=2052\", \"92\u00d795=8740\"),\n    @(\"98\u00d714=1372\", \"92\u00d770=6440\"),\n    @(\"12\u00d742=504\", \"95\u00d764=6080\"),\n    @(\"67\u00d755=3685\", \"92\u00d768=6256\"),\n    @(\"26\u00d776=1976\", \"21\u00d718=378\"),\n    @(\"55\u00d713=715\", \"26\u00d730=780\"),\n    @(\"85\u00d727=2295\", \"73\u00d770=5110\"),\n    @(\"76\u00d736=2736\", \"61\u00d719=1159\"),\n    @(\"47\u00d754=2538\", \"98\u00d719=1862\"),\n    @(\"56\u00d777=4312\", \"36\u00d717=612\"),\n    @(\"57\u00d738=2166\", \"70\u00d799=6930\"),\n    @(\"82\u00d729=2378\", \"73\u00d719=1387\"),\n    @(\"28\u00d776=2128\", \"48\u00d731=1488\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
